$d = $word.ActiveDocument

# Paragraphs (1-indexed, matching Word's Document.Paragraphs collection) whose
# bullet list indent level is increased by one (Increase Indent): they move
# from list level 2 (w:ilvl=1) to list level 3 (w:ilvl=2). These are the
# descriptive detail bullets nested under each field-name bullet in the
# "Database Models" section; the field-name bullets themselves (and the
# model-name / blank bullets) stay put.
$indices = @(
    22, 23, 24, 25,
    27, 28,
    31, 32, 33, 34,
    36, 37,
    39, 40, 41,
    43, 44, 45, 46,
    48, 49,
    52, 53, 54, 55,
    57,
    60, 61, 62, 63,
    65, 66,
    68, 69
)

foreach ($i in $indices) {
    $p = $d.Paragraphs($i)
    $p.Range.ListFormat.ListLevelNumber = 3
}

Write-Output ("Updated " + $indices.Count + " paragraphs to list level 3.")
